$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows for 15, 16 and 17 April 2021 (serial dates 44301-44303)
$newRows = @(
    @{ Row = 227; Date = 44301; B = 0; C = 1; D = 48.07692307692308 },
    @{ Row = 228; Date = 44302; B = 0; C = 1; D = 48.07692307692308 },
    @{ Row = 229; Date = 44303; B = 0; C = 1; D = 48.07692307692308 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D

    # Copy the date-column formatting (border, bold, centered, date number
    # format) from the row above so the new rows match the existing table.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
